$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $find"
    }
    $rng.Text = $replace
}

Replace-Text "These graphics can be used for display in online services, retail establishments, press or other promotional purposes" "Deze afbeeldingen kunnen worden gebruikt voor weergave op online diensten, winkels, pers of voor andere promotionele doeleinden"

Replace-Text ("LOGOS" + [char]0x00A0) "LOGO'S"
Replace-Text ([char]0x00A0 + "ICONS") "ICONEN"

Replace-Text "MERCHANT BUTTONS" "HANDELAAR KNOPPEN"

Replace-Text "SMARTCASH ROADMAP FEATURES" "SMARTCASH ROADMAP FUNCTIES"

Replace-Text "SOCIAL MEDIA GRAPHICS" "SOCIAL MEDIA AFBEELDINGEN"

Replace-Text " GRAPHIC IDENTITY" " GRAFISCHE IDENTITEIT"

Replace-Text "Graphic Identity Guidelines" "Richtlijnen voor grafische edentiteit"
